$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$clothing = @{
    2 = "Trunks,Jumpsuit"
    3 = "Trunks,Dress"
    4 = "Dress,Kaftan"
    5 = "Parka,Jodhpurs"
    6 = "Jodhpurs,Trunks"
    7 = "Blouse,Kaftan"
    8 = "Blazer,Top"
    9 = "Trunks,Jodhpurs"
    10 = "Caftan,Cutoffs"
    11 = "Jumpsuit,Blouse"
    12 = "Jumpsuit,Halter"
    13 = "Trunks,Dress"
    14 = "Trunks,Jumpsuit"
    15 = "Jumpsuit,Dress"
    16 = "Caftan,Trunks"
    17 = "Trunks,Jersey"
    18 = "Jumpsuit,Dress"
    19 = "Parka,Blouse"
    20 = "Jodhpurs,Caftan"
    21 = "Trunks,Jersey"
    22 = "Jodhpurs,Caftan"
    23 = "Jodhpurs,Trunks"
    24 = "Jumpsuit,Blouse"
    25 = "Trunks,Halter"
    26 = "Blouse,Trunks"
    27 = "Trunks,Coverup"
    28 = "Trunks,Cutoffs"
    29 = "Parka,Caftan"
    30 = "Jumpsuit,Trunks"
    31 = "Trunks,Caftan"
    32 = "Trunks,Dress"
    33 = "Caftan,Blazer"
    34 = "Blouse,Top"
    35 = "Trunks,Dress"
    36 = "Trunks,Blouse"
    37 = "Caftan,Trunks"
    38 = "Caftan,Trunks"
    39 = "Trunks,Caftan"
    40 = "Caftan,Trunks"
    41 = "Trunks,Halter"
    42 = "Jumpsuit,Blouse"
    43 = "Trunks,Jumpsuit"
    44 = "Trunks,Dress"
    45 = "Caftan,Trunks"
    46 = "Trunks,Jodhpurs"
    47 = "Kaftan,Blouse"
    48 = "Trunks,Cutoffs"
    49 = "Caftan,Jodhpurs"
    50 = "Jumpsuit,Dress"
    51 = "Jumpsuit,Blouse"
    52 = "Jumpsuit,Halter"
    53 = "Trunks,Jumpsuit"
    54 = "Trunks,Jumpsuit"
    55 = "Trunks,Jumpsuit"
    56 = "Cutoffs,Caftan"
    57 = "Trunks,Caftan"
    58 = "Trunks,Caftan"
    59 = "Sweatpants,Blouse"
    60 = "Trunks,Jodhpurs"
    61 = "Trunks,Sweatpants"
    62 = "Trunks,Kaftan"
    63 = "Trunks,Kaftan"
    64 = "Trunks,Sweatpants"
    65 = "Parka,Caftan"
    66 = "Parka,Caftan"
    67 = "Sweatpants,Parka"
    68 = "Caftan,Jumpsuit"
    69 = "Jumpsuit,Tee"
    70 = "Cutoffs,Blazer"
    71 = "Top,Halter"
    72 = "Trunks,Jumpsuit"
    73 = "Trunks,Jodhpurs"
    74 = "Jumpsuit,Blouse"
    75 = "Blouse,Halter"
    76 = "Caftan,Parka"
    77 = "Trunks,Jodhpurs"
}

foreach ($rowNum in $clothing.Keys) {
    $ws.Cells.Item([int]$rowNum, 7).Value = $clothing[$rowNum]
}
